$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values per the diff
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 7
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 5
